$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($addr, $val)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "63.832.71"
Set-TextCell "E2" "  -3.78%  "
# Row 3
Set-TextCell "D3" "3.342.23"
Set-TextCell "E3" "  -4.44%  "
# Row 4
Set-TextCell "E4" "  +0.30%  "
# Row 5
Set-TextCell "D5" "556.50"
Set-TextCell "E5" "  +0.01%  "
# Row 6
Set-TextCell "D6" "172.04"
Set-TextCell "E6" "  -6.97%  "
# Row 7
Set-TextCell "D7" "0.612"
Set-TextCell "E7" "  -5.49%  "
# Row 8
Set-TextCell "D8" "3.329.14"
Set-TextCell "E8" "  -4.75%  "
# Row 9
Set-TextCell "E9" "  +0.00%  "
# Row 10
Set-TextCell "E10" "  -2.51%  "
# Row 11
Set-TextCell "D11" "0.150"
Set-TextCell "E11" "  -2.37%  "
# Row 12
Set-TextCell "D12" "53.82"
Set-TextCell "E12" "  -1.14%  "
# Row 13
Set-TextCell "D13" "0.0000265"
Set-TextCell "E13" "  -2.27%  "
# Row 14
Set-TextCell "D14" "8.90"
Set-TextCell "E14" "  -4.11%  "
# Row 15
Set-TextCell "D15" "3.886.64"
Set-TextCell "E15" "  -3.94%  "
# Row 16
Set-TextCell "B16" "WrappedEther"
Set-TextCell "C16" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D16" "3.358.32"
Set-TextCell "E16" "  -3.81%  "
# Row 17
Set-TextCell "B17" "TRON"
Set-TextCell "C17" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell "D17" "0.118"
Set-TextCell "E17" "  -3.10%  "
# Row 18
Set-TextCell "D18" "17.71"
Set-TextCell "E18" "  -4.47%  "
# Row 19
Set-TextCell "E19" "  -2.79%  "
# Row 20
Set-TextCell "D20" "63.877.96"
Set-TextCell "E20" "  -3.64%  "
# Row 21
Set-TextCell "D21" "0.972"
Set-TextCell "E21" "  -2.04%  "
# Row 22
Set-TextCell "D22" "403.41"
Set-TextCell "E22" "  -4.50%  "
# Row 23
Set-TextCell "D23" "4.08"
Set-TextCell "E23" "  +0.41%  "
# Row 24
Set-TextCell "D24" "4.29"
Set-TextCell "E24" "  +3.65%  "
# Row 25
Set-TextCell "B25" "InternetComputer(DFINITY)"
Set-TextCell "C25" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D25" "13.30"
Set-TextCell "E25" "  +7.96%  "
# Row 26
Set-TextCell "B26" "Litecoin"
Set-TextCell "C26" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell "D26" "82.76"
Set-TextCell "E26" "  -4.22%  "
# Row 27
Set-TextCell "D27" "10.67"
Set-TextCell "E27" "  -2.60%  "
# Row 28
Set-TextCell "D28" "2.74"
Set-TextCell "E28" "  -5.43%  "
# Row 29
Set-TextCell "D29" "8.72"
Set-TextCell "E29" "  -4.19%  "
# Row 30
Set-TextCell "D30" "29.17"
Set-TextCell "E30" "  -3.32%  "
# Row 31
Set-TextCell "D31" "6.45"
Set-TextCell "E31" "  -2.18%  "
# Row 32
Set-TextCell "D32" "582.88"
Set-TextCell "E32" "  -7.40%  "
# Row 33
Set-TextCell "D33" "11.34"
Set-TextCell "E33" "  -3.22%  "
# Row 34
Set-TextCell "E34" "  -3.58%  "
# Row 35
Set-TextCell "D35" "57.88"
Set-TextCell "E35" "  -3.35%  "
# Row 36
Set-TextCell "E36" "  +0.35%  "
# Row 37
Set-TextCell "D37" "0.999"
Set-TextCell "E37" "  -0.11%  "
# Row 38
Set-TextCell "D38" "35.70"
Set-TextCell "E38" "  -5.03%  "
# Row 39
Set-TextCell "D39" "3.42"
Set-TextCell "E39" "  -0.97%  "
# Row 40
Set-TextCell "D40" "0.0₃0740"
Set-TextCell "E40" "  -8.28%  "
# Row 41
Set-TextCell "D41" "0.369"
Set-TextCell "E41" "  -4.23%  "
# Row 42
Set-TextCell "D42" "3.148.75"
Set-TextCell "E42" "  +1.09%  "
# Row 43
Set-TextCell "D43" "1.00"
Set-TextCell "E43" "  +0.46%  "
# Row 44
Set-TextCell "D44" "2.84"
Set-TextCell "E44" "  -0.17%  "
# Row 45
Set-TextCell "E45" "  -2.40%  "
# Row 46
Set-TextCell "D46" "2.48"
Set-TextCell "E46" "  -4.48%  "
# Row 47
Set-TextCell "D47" "0.0405"
Set-TextCell "E47" "  -2.20%  "
# Row 48
Set-TextCell "D48" "2.62"
Set-TextCell "E48" "  -4.51%  "
# Row 49
Set-TextCell "D49" "0.128"
Set-TextCell "E49" "  -4.12%  "
# Row 50
Set-TextCell "D50" "132.58"
Set-TextCell "E50" "  -4.63%  "
# Row 51
Set-TextCell "D51" "8.07"
Set-TextCell "E51" "  -5.03%  "

Write-Output "applied all changes"
